$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had 24 data rows (rows 2-25). The vessel activity
# mapping used by functions.py changed, which:
#   - re-labels several (scenario, vessel, action) combinations,
#   - re-numbers rows so values shift down within each scenario block,
#   - adds brand-new "Onshore" rows (with no Emissions_factor/Emissions_tons
#     values), and
#   - appends 4 new rows at the end (26-29) for "55 GW" / "Towing Group".
#
# First, extend the sheet with 4 new rows (26-29) and give column A on those
# rows the same header-ish style (s="1": bold, bordered, centered) used by
# the rest of column A, by copying formats from an existing styled cell.
$ws.Range("A25:F25").Copy() | Out-Null
$ws.Range("A26:F29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 2
$ws.Range("A2").Value = "25 GW (SC)"
$ws.Range("B2").Value = "Onshore"
$ws.Range("C2").Value = "Idle at port"
$ws.Range("D2").Value = 24635.25
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()

# Row 3
$ws.Range("A3").Value = "25 GW (SC)"
$ws.Range("B3").Value = "Onshore"
$ws.Range("C3").Value = "None"
$ws.Range("D3").Value = 424809
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()

# Row 4
$ws.Range("A4").Value = "25 GW (SC)"
$ws.Range("B4").Value = "Onshore"
$ws.Range("C4").Value = "Transit"
$ws.Range("D4").Value = 11088
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()

# Row 5
$ws.Range("A5").Value = "25 GW (SC)"
$ws.Range("B5").Value = "Towing Group"
$ws.Range("C5").Value = "Idle at port"
$ws.Range("D5").Value = 678265.29262013
$ws.Range("E5").Value = 0.407310331
$ws.Range("F5").Value = 276264.460842917

# Row 6
$ws.Range("A6").Value = "25 GW (SC)"
$ws.Range("B6").Value = "Towing Group"
$ws.Range("C6").Value = "Idle at sea"
$ws.Range("D6").Value = 108108
$ws.Range("E6").Value = 0.203655165
$ws.Range("F6").Value = 22016.75257782

# Row 7
$ws.Range("A7").Value = "25 GW (SC)"
$ws.Range("B7").Value = "Towing Group"
$ws.Range("C7").Value = "Maneuvering"
$ws.Range("D7").Value = 8316
$ws.Range("E7").Value = 2.698858249
$ws.Range("F7").Value = 22443.705198684

# Row 8
$ws.Range("A8").Value = "25 GW (SC)"
$ws.Range("B8").Value = "Towing Group"
$ws.Range("C8").Value = "Transit"
$ws.Range("D8").Value = 198897.0324545454
$ws.Range("E8").Value = 4.753093345
$ws.Range("F8").Value = 945376.1612999489

# Row 9
$ws.Range("A9").Value = "25 GW (CC)"
$ws.Range("B9").Value = "Onshore"
$ws.Range("C9").Value = "Idle at port"
$ws.Range("D9").Value = 25368.5
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()

# Row 10
$ws.Range("A10").Value = "25 GW (CC)"
$ws.Range("B10").Value = "Onshore"
$ws.Range("C10").Value = "None"
$ws.Range("D10").Value = 424809
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()

# Row 11
$ws.Range("A11").Value = "25 GW (CC)"
$ws.Range("B11").Value = "Onshore"
$ws.Range("C11").Value = "Transit"
$ws.Range("D11").Value = 11088
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()

# Row 12
$ws.Range("A12").Value = "25 GW (CC)"
$ws.Range("B12").Value = "Towing Group"
$ws.Range("C12").Value = "Idle at port"
$ws.Range("D12").Value = 672579.5244204547
$ws.Range("E12").Value = 0.407310331
$ws.Range("F12").Value = 273948.588715518

# Row 13
$ws.Range("A13").Value = "25 GW (CC)"
$ws.Range("B13").Value = "Towing Group"
$ws.Range("C13").Value = "Idle at sea"
$ws.Range("D13").Value = 108108
$ws.Range("E13").Value = 0.203655165
$ws.Range("F13").Value = 22016.75257782

# Row 14
$ws.Range("A14").Value = "25 GW (CC)"
$ws.Range("B14").Value = "Towing Group"
$ws.Range("C14").Value = "Maneuvering"
$ws.Range("D14").Value = 8316
$ws.Range("E14").Value = 2.698858249
$ws.Range("F14").Value = 22443.705198684

# Row 15
$ws.Range("A15").Value = "25 GW (CC)"
$ws.Range("B15").Value = "Towing Group"
$ws.Range("C15").Value = "Transit"
$ws.Range("D15").Value = 72408.375
$ws.Range("E15").Value = 4.753093345
$ws.Range("F15").Value = 344163.7653347644

# Row 16
$ws.Range("A16").Value = "35 GW"
$ws.Range("B16").Value = "Onshore"
$ws.Range("C16").Value = "Idle at port"
$ws.Range("D16").Value = 34547.625
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()

# Row 17
$ws.Range("A17").Value = "35 GW"
$ws.Range("B17").Value = "Onshore"
$ws.Range("C17").Value = "None"
$ws.Range("D17").Value = 593690.5
$ws.Range("E17").ClearContents()
$ws.Range("F17").ClearContents()

# Row 18
$ws.Range("A18").Value = "35 GW"
$ws.Range("B18").Value = "Onshore"
$ws.Range("C18").Value = "Transit"
$ws.Range("D18").Value = 15496
$ws.Range("E18").ClearContents()
$ws.Range("F18").ClearContents()

# Row 19
$ws.Range("A19").Value = "35 GW"
$ws.Range("B19").Value = "Towing Group"
$ws.Range("C19").Value = "Idle at port"
$ws.Range("D19").Value = 929730.2761655845
$ws.Range("E19").Value = 0.407310331
$ws.Range("F19").Value = 378688.7465257256

# Row 20
$ws.Range("A20").Value = "35 GW"
$ws.Range("B20").Value = "Towing Group"
$ws.Range("C20").Value = "Idle at sea"
$ws.Range("D20").Value = 151086
$ws.Range("E20").Value = 0.203655165
$ws.Range("F20").Value = 30769.44425919

# Row 21
$ws.Range("A21").Value = "35 GW"
$ws.Range("B21").Value = "Towing Group"
$ws.Range("C21").Value = "Maneuvering"
$ws.Range("D21").Value = 11622
$ws.Range("E21").Value = 2.698858249
$ws.Range("F21").Value = 31366.130569878

# Row 22
$ws.Range("A22").Value = "35 GW"
$ws.Range("B22").Value = "Towing Group"
$ws.Range("C22").Value = "Transit"
$ws.Range("D22").Value = 244802.1761883117
$ws.Range("E22").Value = 4.753093345
$ws.Range("F22").Value = 1163567.594482182

# Row 23
$ws.Range("A23").Value = "55 GW"
$ws.Range("B23").Value = "Onshore"
$ws.Range("C23").Value = "Idle at port"
$ws.Range("D23").Value = 52839
$ws.Range("E23").ClearContents()
$ws.Range("F23").ClearContents()

# Row 24
$ws.Range("A24").Value = "55 GW"
$ws.Range("B24").Value = "Onshore"
$ws.Range("C24").Value = "None"
$ws.Range("D24").Value = 900190.5
$ws.Range("E24").ClearContents()
$ws.Range("F24").ClearContents()

# Row 25
$ws.Range("A25").Value = "55 GW"
$ws.Range("B25").Value = "Onshore"
$ws.Range("C25").Value = "Transit"
$ws.Range("D25").Value = 23496
$ws.Range("E25").ClearContents()
$ws.Range("F25").ClearContents()

# Row 26
$ws.Range("A26").Value = "55 GW"
$ws.Range("B26").Value = "Towing Group"
$ws.Range("C26").Value = "Idle at port"
$ws.Range("D26").Value = 1523056.469892857
$ws.Range("E26").Value = 0.407310331
$ws.Range("F26").Value = 620356.6348837512

# Row 27
$ws.Range("A27").Value = "55 GW"
$ws.Range("B27").Value = "Towing Group"
$ws.Range("C27").Value = "Idle at sea"
$ws.Range("D27").Value = 229086
$ws.Range("E27").Value = 0.203655165
$ws.Range("F27").Value = 46654.54712919

# Row 28
$ws.Range("A28").Value = "55 GW"
$ws.Range("B28").Value = "Towing Group"
$ws.Range("C28").Value = "Maneuvering"
$ws.Range("D28").Value = 17622
$ws.Range("E28").Value = 2.698858249
$ws.Range("F28").Value = 47559.280063878

# Row 29
$ws.Range("A29").Value = "55 GW"
$ws.Range("B29").Value = "Towing Group"
$ws.Range("C29").Value = "Transit"
$ws.Range("D29").Value = 536067.7882012987
$ws.Range("E29").Value = 4.753093345
$ws.Range("F29").Value = 2547980.236568462

